# Daily attendance processing - 2026-01-18 09:59:59
#
# Swap the order of "System" and the recorder's email address in every
# "Recorded By" cell (column G) that currently reads
# "System, dnasr281@gmail.com", turning it into
# "dnasr281@gmail.com, System".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$lastRow = $ws.UsedRange.Rows.Count
$updated = 0

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
        $updated = $updated + 1
    }
}

Write-Host "Updated $updated 'Recorded By' cells from '$oldValue' to '$newValue'"
